# Generate Report for Handoff
#
# Updates the localization-status report:
#  - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#    on every sheet (Overview!E2/F2, zh-cn!C2, de-de!C2)
#  - The "Latest Handoff Datetime" timestamps are refreshed to reflect the
#    new handoff generation time (Overview!G2, zh-cn!H2, de-de!H2)
#  - The Status / zh-cn / de-de columns are narrowed to fit the new,
#    shorter "Ready for handoff" text

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$newWidth = 17.2159881591797

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-30 11:08:13"
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-30 11:08:05"
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-30 11:08:13"
$dede.Columns.Item(3).ColumnWidth = $newWidth
